$descriptionUpdates = @(
    @{ Row = 8; Text = 'I wanted the ability to change the grid size while the game was running instead of having to rerun it the entire time. I first began by making a simple UI that had 2 fields where a person could input their desired columns and rows, a button which would generate the grid and a button which would start the generation of the maze. I did this with the Unity UI Toolkit. I already have some experience with it and, since it uses CSS and Flexbox, the UI could easily be turned responsive. ' }
    @{ Row = 9; Text = 'Buttons and fields get assigned to their respective variables and click events if needed. To make clients have the ability to call the needed functions in the maze generation script, I, for now, have added a reference to the UI script to quickly test if the buttons work. The trigger for starting the maze generation corountine was moved from start to their own function so that it only starts when the user has filled in their desired size.' }
    @{ Row = 10; Text = 'Had an issue where sometimes the index of the 2D array would be out of bounds when newly generated, breaking the whole thing. This was fixed by, when removing the grid, setting the 2d array to null and reassigning it when the maze would be generated again.' }
    @{ Row = 11; Text = 'One of the conditions given for the assessment states that the player must be able to see the maze in its entirety. Since the maze could be as big as 250x250 or as small as 10x10, the maze could or, be outside the camera view window or, be too small to enjoyably see the maze. I therefore decided to give the player the ability to move the camera up and down, basically making them be able to zoom in and out on a maze. This was done via adding the necessary force to the rigid body of the player.
' }
    @{ Row = 12; Text = 'Expanding on the previous idea, I thought it would be fun to see the maze you generated from different angles, not just from above. I therefore implemented simple player movement to the game using this tutorial: https://youtu.be/f473C43s8nE?si=mpyzDA6Pd4wtbr_A, the Unity Input System and reused code from an old project of mine. The input system registers the pressed buttons of the player and calls upon the function connected to it via unity events. These button presses then get converted to rigid body force which pushes the player to the desired direction. I kept in the movement I made for the up and down movement of the player so that they could fly around freely.' }
    @{ Row = 13; Text = 'Implemented the ability for the user to switch between UI mode and Moving mode via pressing TAB. When the player hits TAB the camera stops following the cursor and the cursor becomes visible again, giving the player the ability to configure the maze. If the player hits TAB again, the camera locks onto the movement of the cursor again and the cursor becomes invisible again.' }
    @{ Row = 14; Text = 'When the player is moving, the buttons and input field would get selected and cause the player to accidentally input a bunch of button presses in the input field. I had some difficulty fixing this, but eventually I noticed that unity hadn''t generated an event system in the hierarchy yet. When I did notice and made one, I instantly fixed it by unassigning moving from the unity events list.' }
    @{ Row = 15; Text = 'For the bigger mazes, the generation could take a few seconds to be ready. To show the player that it was still running and how far it already was, I added a progress bar which fills up till the maze is completed. I did this using the Unity UI Toolkit and assigning the values to the progressbar when the maze started generating.' }
    @{ Row = 16; Text = 'As an extra feature, I thought it would be fun for players to have the ability to shoot down the walls of the generated maze so they could make their own path. I did this by having the player shoot out a raycast when they clicked the left mouse button. The raycast would then detect hits with objects that have the MazeWall Layer added and promptly delete them from the scene.' }
    @{ Row = 17; Text = 'After testing, I noticed that grids larger than 60x60 would begin to tank performance by quite a bit. I first started by adding batching to the cells once they were initialised, which combined game objects with the same materials. But this didn''t help as much as I would''ve liked so I made the decision to have the grid first generate invisible. The algorithm would then, when activated, turn on the walls that were needed. Because of the reduced number of game object on the screen FPS doubled compared to when the grid was visible the entire time.' }
    @{ Row = 18; Text = 'Fixed an issue with the Maze algorithm causing a stackoverflow after a certain size. This was caused by the maze algorithm starting too many coroutines after each other, which would fill up the limit. It was eventually fixed after I found out I had placed an extra yield return before the while do loop. Removing this extra yield statement fixed the issue.' }
    @{ Row = 19; Text = 'Implemented the UI to be responsive to the screen size. Due to how UI Toolkit works, I only had to change around some absolute sizes to percentages to make it work without issues.' }
    @{ Row = 20; Text = 'To lessen the dependency of some scripts, I reconfigured some references to scripts to work via unity events. The workings of these events are similar to how the Input System works. One script invokes an event with a custom-named scriptable object. The event listener that''s attached to the script that must do something based upon that event picks up the call and compares it with their custom-named scriptable object. If they are the same, the listener activates the assigned function in that script. ' }
    @{ Row = 21; Text = 'The code for deactivating different walls had a lot of duplicate code that looked quite bulky. I therefore decided to add all the different wall types (Up, Down, Left, Right) into a list in the cells and merge all wall removal functions together into 1. This new WallRemove function would get an enum int from the maze generation script and compare it with the Wall type list. This wall would then be set to not be enabled again via a bool. After this, the function will go through the wall list and enable each wall which isn''t set to not be enabled again. This made the code a lot more readable and flexible to expand if needed.' }
    @{ Row = 22; Text = 'While testing around and destroying walls, I sometimes accidentally shot the wrong wall. This was mildly annoying as it would mean I would need to generate the maze again to sort of fix that mistake. I then decided that I should also add the ability to add walls back again after you shot them. To achieve this, I had to change the logic for the wall destruction to, instead of removing it from the scene, disable the mesh render and turn on the isTrigger variable in the box collider of the wall (Otherwise the ray wouldn''t be able to detect the collision any more). The logic for this feature is basically the same as the logic for wall destruction but then with flipped variables. ' }
    @{ Row = 23; Text = 'During testing, I noticed that the game would freeze every time you (re)generated a large maze. This was caused by the grid used being initialised and destroyed each time a maze was generated. To stop the game from freezing each time you generate something larger than a 50x50, I decided to initialise a 250x250 (invisible) grid at the start of the game. If a player decides to regenerate a maze, all the walls first get reset to being invisible again before the algorithm begins generating again. The gird being initialised at the beginning had the added benefit of the player being able to build walls outside their generated maze (as long as it''s smaller than 250x250).' }
    @{ Row = 24; Text = 'Crafted up a small low-poly world with slopes surrounding the play area and a small shrine where the player spawns to give the game a bit more than just a grey plane' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($update in $descriptionUpdates) {
    $ws.Cells.Item($update.Row, 4).Value2 = $update.Text
}

# Row 11's description text grew longer, so Excel auto-expanded its row height to fit
$ws.Rows.Item(11).RowHeight = 122.25

# Update the view state: scroll position and active cell selection
$ws.Activate()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 19
    $win.ScrollColumn = 1
} catch {}
$ws.Range("G25").Select()
